$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add a new comment on E12 ---
$ws.Range("E12").AddComment("Not implemented in a servlet. Otherwise successfully extracting filtered requests from JSON data.")

# --- Update task #8 (row 12) and task #9 (row 13): mark as Started & assigned to Svetoslav ---
$ws.Range("A12:E13").Interior.Color = 65535

$ws.Range("C12").Value = "Svetoslav"
$ws.Range("D12").Value = "Started"

$ws.Range("C13").Value = "Svetoslav"
$ws.Range("D13").Value = "Started"

# --- Update the sheet view: scroll position and active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F14").Select()
